$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7791272401809692
$ws.Range("B1").Value = 1.271139144897461
$ws.Range("C1").Value = 5.13085412979126
$ws.Range("D1").Value = 1.607800602912903
$ws.Range("E1").Value = 0.9239203929901123
